# Update PLC data 2025-10-13 14:10:36
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7202
$ws.Range("C3").Value = 174548
$ws.Range("C4").Value = 164541
$ws.Range("C7").Value = 5.73
$ws.Range("C8").Value = 64.37
